# Lecture10.pptx - "small changes to lecture #10"
#
# 1) The cached text of every automatic "date" field (the
#    datetimeFigureOut field living on the slide master, all 11
#    slide layouts, and the notes master) is refreshed from
#    2/20/2018 to 2/24/2018.
# 2) A straight arrow connector on slide 35 is nudged to the right
#    (its x-offset moves from 2895600 EMU to 3200400 EMU, i.e. from
#    228pt to 252pt; the y-offset is untouched).

$p = $ppt.ActivePresentation

$oldDate = "2/20/2018"
$newDate = "2/24/2018"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's Date Placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's Date Placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master's Date Placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# Move the straight arrow connector on slide 35.
$slide = $p.Slides.Item(35)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Straight Arrow Connector 11") {
        $shp.Left = 252.0
    }
}
